$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (B1:E1) changed from 1,2,3,4 to 15,16,15,16
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (meanEMG legmaxROM) updated values for B2:E2
$ws.Range("B2").Value = 107.30999307531243
$ws.Range("C2").Value = 105.06374212778209
$ws.Range("D2").Value = 106.18689007183453
$ws.Range("E2").Value = 106.94149088775504

# Row 3 updated values for B3:E3
$ws.Range("B3").Value = 106.28805611779843
$ws.Range("C3").Value = 104.45507621814201
$ws.Range("D3").Value = 106.06479608680182
$ws.Range("E3").Value = 107.26461657240205

# Selection now only covers B1:E3 instead of B1:AY3
$ws.Range("B1:E3").Select()
